$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = 0.01636035618807323
$ws.Range("C2").Value = 0.7542381384319475
$ws.Range("D2").Value = 1.208071944140944
$ws.Range("E2").Value = 1.099123261577583
$ws.Range("F2").Value = 1.111419983725625
$ws.Range("G2").Value = 45

# Row 3 (Q0)
$ws.Range("B3").Value = 0.1050767557480483
$ws.Range("C3").Value = 1.17191445325214
$ws.Range("D3").Value = 3.506314857596029
$ws.Range("E3").Value = 1.872515649492957
$ws.Range("F3").Value = 1.876375943129116

# Row 4 (Q1)
$ws.Range("B4").Value = 0.1730876980324934
$ws.Range("C4").Value = 1.326443179105317
$ws.Range("D4").Value = 7.8735929385556
$ws.Range("E4").Value = 2.805992326888226
$ws.Range("F4").Value = 2.822110003756866
$ws.Range("G4").Value = 66
